$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - copy H1's formatting (bold/border/center/top) onto the new
# header cells, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF), rows 2-10
$values = @{
    2  = @(5, 7)
    3  = @(4, 5)
    4  = @(1, 1)
    5  = @(5, 6)
    6  = @(4, 4)
    7  = @(2, 5)
    8  = @(8, 9)
    9  = @(8, 8)
    10 = @(6, 6)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
